# Quotient_Ext.xlsx edit: prefix the "Variable Name" column (B) for rows 5-46
# with "Quotient_" (DX explanations update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 5; $r -le 46; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value = "Quotient_" + $cell.Value2
}

# Reset the active selection to A3 (matches the saved view state).
$ws.Range("A3").Select()
